$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.947.96"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "2.969.10"
$ws.Range("E3").Value = "  -2.35%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.57"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.42"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "2.966.09"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.20"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.148"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.17"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "3.447.03"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "60.953.24"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "2.965.06"
$ws.Range("E18").Value = "  -2.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "455.31"
$ws.Range("E20").Value = "  -3.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.08"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.669"
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.79"
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.65"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.69"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.62"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.60"
$ws.Range("E28").Value = "  -5.48%  "
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.25"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.82"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.11"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "54.60"
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.23"
$ws.Range("E34").Value = "  -5.46%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("E35").Value = "  +3.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.73"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "453.44"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").Value = "3.129.96"
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0772"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0379"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.117"
$ws.Range("E41").Value = "  +5.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.93"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.36"
$ws.Range("E44").Value = "  -6.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.242"
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.18"
$ws.Range("E46").Value = "  +5.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.88"
$ws.Range("E47").Value = "  +3.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.107"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.91"
$ws.Range("E49").Value = "  -4.14%  "
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").Value = "0.0₃0500"
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("E51").Value = "  +5.64%  "
